# 电信主要通信能力 - refresh data from 2000-2020 to 2010-2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer part of the (shorter) refreshed series.
# Before: rows 2..22 held 2000..2020 (21 years). After: rows 2..14 hold
# 2010..2022 (13 years), so the old rows 15..22 disappear entirely.
$ws.Rows("15:22").Delete()

# New data for rows 2..14 (years 2010..2022).
$data = @(
    @("2010年", 9962466.5,    1641.5,    46537.3,     150284.9,    81.8133),
    @("2011年", 12119302.9,   1602.3,    43428.4,     171636,      84.23408999999999),
    @("2012年", 14793300.43,  1579.7,    43749.3079,  184023.824,  86.817525),
    @("2013年", 17453709.2,   1280.5,    41089.3,     196557.3,    89.00182),
    @("2014年", 20612529.22,  982.9,     40517.14,    205024.92,   92.83982),
    @("2015年", 24863348.24,  811.1,     26446.45,    218149.97,   96.5283),
    @("2016年", 30420755.06,  681.1,     22441.59,    218540.03,   99.409161),
    @("2017年", 37801073.37,  603.5297,  18398.7466,  242185.7823, 104.499781),
    @("2018年", 43167888.41,  392.4342,  11440.4,     259453.14,   99.413014),
    @("2019年", 47412442.31,  119.3694,  7189.7423,   272523.75,   108.493705),
    @("2020年", 51692051.4,   74.0454,   6923.84,     274567.13,   111.792297)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Row 13: 2021年 -- the long-distance switch columns (C, D) are no longer
# reported (blank), while B, E and F keep being populated.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 54808232.79
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 275690.78
$ws.Range("F13").Value = 112.083653

# Row 14: 2022年 -- only the cumulative total (B) is available so far.
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 59580000
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
